$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.016736319288611412
$ws.Range("C2").Value = 0.00830896571278572
$ws.Range("D2").Value = 0.005099698901176453
$ws.Range("E2").Value = 0.00380017701536417
$ws.Range("F2").Value = 0.0001602680131327361
$ws.Range("I2").Value = 1.2575732469558716
$ws.Range("J2").Value = 0.12731222808361053
$ws.Range("K2").Value = 1.4204224348068237
